$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table of scalar parameter values grows from one row to nine. Row 2
# (originally "param_pv2_area" / 200) now reads "param_pv1_area" and keeps
# the value 200, while a new row 3 takes over the "param_pv2_area" label
# with its own value (250); rows 4-10 are brand new parameters appended
# below.
$rows = @(
    @{ Name = "param_pv1_area"; Value = 200 },
    @{ Name = "param_pv2_area"; Value = 250 },
    @{ Name = "param_bat1_E_max_initial"; Value = 1000 },
    @{ Name = "param_bat2_E_max_initial"; Value = 100 },
    @{ Name = "param_solar_th2_area"; Value = 300 },
    @{ Name = "param_pvt1_area"; Value = 120 },
    @{ Name = "param_pvt2_area"; Value = 150 },
    @{ Name = "param_Q_gas_boiler1_max"; Value = 20 },
    @{ Name = "param_Q_gas_boiler2_max"; Value = 100 }
)

# Copy the existing label cell's formatting (bold, centered, bordered style)
# down onto all the new label rows before filling in values, so every new
# cell picks up the same cell style the sheet already uses for A2/B1.
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Name
    $ws.Cells.Item($r, 2).Value = $row.Value
    $r++
}
